$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -19.21856323619027
$ws.Cells.Item(2, 3).Value = 1.813283995637452
$ws.Cells.Item(2, 4).Value = -19.21856323619027
$ws.Cells.Item(2, 5).Value = -19.21856323619027
$ws.Cells.Item(2, 6).Value = -19.21856323619027
$ws.Cells.Item(2, 7).Value = -19.21856323619027
$ws.Cells.Item(2, 8).Value = -19.21856323619027
$ws.Cells.Item(2, 9).Value = -19.21856323619027
$ws.Cells.Item(2, 10).Value = -19.21856323619027
$ws.Cells.Item(2, 11).Value = -19.21856323619027
$ws.Cells.Item(3, 2).Value = -19.21856323619027
$ws.Cells.Item(3, 3).Value = -19.21856323619027
$ws.Cells.Item(3, 4).Value = -19.21856323619027
$ws.Cells.Item(3, 5).Value = -19.21856323619027
$ws.Cells.Item(3, 6).Value = -19.21856323619027
$ws.Cells.Item(3, 7).Value = -19.21856323619027
$ws.Cells.Item(3, 8).Value = -19.21856323619027
$ws.Cells.Item(3, 9).Value = 2.669359789531192
$ws.Cells.Item(3, 10).Value = -19.21856323619027
$ws.Cells.Item(3, 11).Value = -19.21856323619027
$ws.Cells.Item(4, 2).Value = -19.21856323619027
$ws.Cells.Item(4, 3).Value = 2.155806078024904
$ws.Cells.Item(4, 4).Value = 2.169812007936245
$ws.Cells.Item(4, 5).Value = -19.21856323619027
$ws.Cells.Item(4, 6).Value = 3.501528392113449
$ws.Cells.Item(4, 7).Value = -19.21856323619027
$ws.Cells.Item(4, 8).Value = 1.688915007627182
$ws.Cells.Item(4, 9).Value = -19.21856323619027
$ws.Cells.Item(4, 10).Value = 2.162730974150718
$ws.Cells.Item(4, 11).Value = -19.21856323619027
$ws.Cells.Item(5, 2).Value = -19.21856323619027
$ws.Cells.Item(5, 3).Value = 1.868884713308972
$ws.Cells.Item(5, 4).Value = -19.21856323619027
$ws.Cells.Item(5, 5).Value = -19.21856323619027
$ws.Cells.Item(5, 6).Value = -19.21856323619027
$ws.Cells.Item(5, 7).Value = 2.835585623758208
$ws.Cells.Item(5, 8).Value = -19.21856323619027
$ws.Cells.Item(5, 9).Value = -19.21856323619027
$ws.Cells.Item(5, 10).Value = -19.21856323619027
$ws.Cells.Item(5, 11).Value = -19.21856323619027
$ws.Cells.Item(6, 2).Value = -19.21856323619027
$ws.Cells.Item(6, 3).Value = -19.21856323619027
$ws.Cells.Item(6, 4).Value = -19.21856323619027
$ws.Cells.Item(6, 5).Value = -19.21856323619027
$ws.Cells.Item(6, 6).Value = -19.21856323619027
$ws.Cells.Item(6, 7).Value = -19.21856323619027
$ws.Cells.Item(6, 8).Value = -19.21856323619027
$ws.Cells.Item(6, 9).Value = -19.21856323619027
$ws.Cells.Item(6, 10).Value = -19.21856323619027
$ws.Cells.Item(6, 11).Value = -19.21856323619027
$ws.Cells.Item(7, 2).Value = -19.21856323619027
$ws.Cells.Item(7, 3).Value = -19.21856323619027
$ws.Cells.Item(7, 4).Value = -19.21856323619027
$ws.Cells.Item(7, 5).Value = -19.21856323619027
$ws.Cells.Item(7, 6).Value = -19.21856323619027
$ws.Cells.Item(7, 7).Value = -19.21856323619027
$ws.Cells.Item(7, 8).Value = -19.21856323619027
$ws.Cells.Item(7, 9).Value = -19.21856323619027
$ws.Cells.Item(7, 10).Value = -19.21856323619027
$ws.Cells.Item(7, 11).Value = -19.21856323619027
$ws.Cells.Item(8, 2).Value = -19.21856323619027
$ws.Cells.Item(8, 3).Value = -19.21856323619027
$ws.Cells.Item(8, 4).Value = -19.21856323619027
$ws.Cells.Item(8, 5).Value = 1.900596374048489
$ws.Cells.Item(8, 6).Value = -19.21856323619027
$ws.Cells.Item(8, 7).Value = -19.21856323619027
$ws.Cells.Item(8, 8).Value = -19.21856323619027
$ws.Cells.Item(8, 9).Value = -19.21856323619027
$ws.Cells.Item(8, 10).Value = -19.21856323619027
$ws.Cells.Item(8, 11).Value = -19.21856323619027
$ws.Cells.Item(9, 2).Value = 4.321925848242253
$ws.Cells.Item(9, 3).Value = -19.21856323619027
$ws.Cells.Item(9, 4).Value = -19.21856323619027
$ws.Cells.Item(9, 5).Value = -19.21856323619027
$ws.Cells.Item(9, 6).Value = -19.21856323619027
$ws.Cells.Item(9, 7).Value = -19.21856323619027
$ws.Cells.Item(9, 8).Value = -19.21856323619027
$ws.Cells.Item(9, 9).Value = -19.21856323619027
$ws.Cells.Item(9, 10).Value = -19.21856323619027
$ws.Cells.Item(9, 11).Value = -19.21856323619027
$ws.Cells.Item(10, 2).Value = -19.21856323619027
$ws.Cells.Item(10, 3).Value = -19.21856323619027
$ws.Cells.Item(10, 4).Value = -19.21856323619027
$ws.Cells.Item(10, 5).Value = -19.21856323619027
$ws.Cells.Item(10, 6).Value = -19.21856323619027
$ws.Cells.Item(10, 7).Value = -19.21856323619027
$ws.Cells.Item(10, 8).Value = -19.21856323619027
$ws.Cells.Item(10, 9).Value = 1.211686136318545
$ws.Cells.Item(10, 10).Value = -19.21856323619027
$ws.Cells.Item(10, 11).Value = 1.983977331631268
$ws.Cells.Item(11, 2).Value = -19.21856323619027
$ws.Cells.Item(11, 3).Value = -19.21856323619027
$ws.Cells.Item(11, 4).Value = -19.21856323619027
$ws.Cells.Item(11, 5).Value = 2.87343539921395
$ws.Cells.Item(11, 6).Value = -19.21856323619027
$ws.Cells.Item(11, 7).Value = 2.872731707088671
$ws.Cells.Item(11, 8).Value = -19.21856323619027
$ws.Cells.Item(11, 9).Value = -19.21856323619027
$ws.Cells.Item(11, 10).Value = -19.21856323619027
$ws.Cells.Item(11, 11).Value = 1.845016831278732
$ws.Cells.Item(12, 2).Value = -19.21856323619027
$ws.Cells.Item(12, 3).Value = -19.21856323619027
$ws.Cells.Item(12, 4).Value = -19.21856323619027
$ws.Cells.Item(12, 5).Value = -19.21856323619027
$ws.Cells.Item(12, 6).Value = -19.21856323619027
$ws.Cells.Item(12, 7).Value = -19.21856323619027
$ws.Cells.Item(12, 8).Value = -19.21856323619027
$ws.Cells.Item(12, 9).Value = -19.21856323619027
$ws.Cells.Item(12, 10).Value = -19.21856323619027
$ws.Cells.Item(12, 11).Value = -19.21856323619027
$ws.Cells.Item(13, 2).Value = -19.21856323619027
$ws.Cells.Item(13, 3).Value = -19.21856323619027
$ws.Cells.Item(13, 4).Value = -19.21856323619027
$ws.Cells.Item(13, 5).Value = 2.446099748280209
$ws.Cells.Item(13, 6).Value = -19.21856323619027
$ws.Cells.Item(13, 7).Value = -19.21856323619027
$ws.Cells.Item(13, 8).Value = -19.21856323619027
$ws.Cells.Item(13, 9).Value = -19.21856323619027
$ws.Cells.Item(13, 10).Value = 1.963070096068121
$ws.Cells.Item(13, 11).Value = 1.927230253074014
$ws.Cells.Item(14, 2).Value = -19.21856323619027
$ws.Cells.Item(14, 3).Value = -19.21856323619027
$ws.Cells.Item(14, 4).Value = 1.115357512751011
$ws.Cells.Item(14, 5).Value = -19.21856323619027
$ws.Cells.Item(14, 6).Value = -19.21856323619027
$ws.Cells.Item(14, 7).Value = -19.21856323619027
$ws.Cells.Item(14, 8).Value = -19.21856323619027
$ws.Cells.Item(14, 9).Value = -19.21856323619027
$ws.Cells.Item(14, 10).Value = -19.21856323619027
$ws.Cells.Item(14, 11).Value = 2.147247100242697
$ws.Cells.Item(15, 2).Value = -19.21856323619027
$ws.Cells.Item(15, 3).Value = -19.21856323619027
$ws.Cells.Item(15, 4).Value = 1.122674326186674
$ws.Cells.Item(15, 5).Value = -19.21856323619027
$ws.Cells.Item(15, 6).Value = -19.21856323619027
$ws.Cells.Item(15, 7).Value = -19.21856323619027
$ws.Cells.Item(15, 8).Value = -19.21856323619027
$ws.Cells.Item(15, 9).Value = -19.21856323619027
$ws.Cells.Item(15, 10).Value = -19.21856323619027
$ws.Cells.Item(15, 11).Value = -19.21856323619027
$ws.Cells.Item(16, 2).Value = -19.21856323619027
$ws.Cells.Item(16, 3).Value = -19.21856323619027
$ws.Cells.Item(16, 4).Value = -19.21856323619027
$ws.Cells.Item(16, 5).Value = -19.21856323619027
$ws.Cells.Item(16, 6).Value = -19.21856323619027
$ws.Cells.Item(16, 7).Value = -19.21856323619027
$ws.Cells.Item(16, 8).Value = -19.21856323619027
$ws.Cells.Item(16, 9).Value = -19.21856323619027
$ws.Cells.Item(16, 10).Value = 2.086470359981207
$ws.Cells.Item(16, 11).Value = -19.21856323619027
$ws.Cells.Item(17, 2).Value = -19.21856323619027
$ws.Cells.Item(17, 3).Value = 2.151224184803927
$ws.Cells.Item(17, 4).Value = 2.483970081240615
$ws.Cells.Item(17, 5).Value = -19.21856323619027
$ws.Cells.Item(17, 6).Value = -19.21856323619027
$ws.Cells.Item(17, 7).Value = -19.21856323619027
$ws.Cells.Item(17, 8).Value = 1.370353568923766
$ws.Cells.Item(17, 9).Value = 2.033606115269679
$ws.Cells.Item(17, 10).Value = 2.098547363838816
$ws.Cells.Item(17, 11).Value = -19.21856323619027
$ws.Cells.Item(18, 2).Value = -19.21856323619027
$ws.Cells.Item(18, 3).Value = -19.21856323619027
$ws.Cells.Item(18, 4).Value = -19.21856323619027
$ws.Cells.Item(18, 5).Value = -19.21856323619027
$ws.Cells.Item(18, 6).Value = -19.21856323619027
$ws.Cells.Item(18, 7).Value = -19.21856323619027
$ws.Cells.Item(18, 8).Value = 1.632140200734911
$ws.Cells.Item(18, 9).Value = 1.313390231190266
$ws.Cells.Item(18, 10).Value = 1.629283999476605
$ws.Cells.Item(18, 11).Value = -19.21856323619027
$ws.Cells.Item(19, 2).Value = -19.21856323619027
$ws.Cells.Item(19, 3).Value = -19.21856323619027
$ws.Cells.Item(19, 4).Value = 1.574310444988922
$ws.Cells.Item(19, 5).Value = -19.21856323619027
$ws.Cells.Item(19, 6).Value = -19.21856323619027
$ws.Cells.Item(19, 7).Value = -19.21856323619027
$ws.Cells.Item(19, 8).Value = 1.557006055409883
$ws.Cells.Item(19, 9).Value = 1.537456308362777
$ws.Cells.Item(19, 10).Value = -19.21856323619027
$ws.Cells.Item(19, 11).Value = -19.21856323619027
$ws.Cells.Item(20, 2).Value = -19.21856323619027
$ws.Cells.Item(20, 3).Value = 0.8630723297194925
$ws.Cells.Item(20, 4).Value = 1.369639902905352
$ws.Cells.Item(20, 5).Value = -19.21856323619027
$ws.Cells.Item(20, 6).Value = 3.116742332292592
$ws.Cells.Item(20, 7).Value = -19.21856323619027
$ws.Cells.Item(20, 8).Value = 1.905717765328732
$ws.Cells.Item(20, 9).Value = 0.8797734724579633
$ws.Cells.Item(20, 10).Value = -19.21856323619027
$ws.Cells.Item(20, 11).Value = 2.076737676274721
$ws.Cells.Item(21, 2).Value = -19.21856323619027
$ws.Cells.Item(21, 3).Value = 1.081294287057927
$ws.Cells.Item(21, 4).Value = -19.21856323619027
$ws.Cells.Item(21, 5).Value = 1.802734649196842
$ws.Cells.Item(21, 6).Value = -19.21856323619027
$ws.Cells.Item(21, 7).Value = 2.469140921118977
$ws.Cells.Item(21, 8).Value = 2.137960197284193
$ws.Cells.Item(21, 9).Value = -19.21856323619027
$ws.Cells.Item(21, 10).Value = -19.21856323619027
$ws.Cells.Item(21, 11).Value = -19.21856323619027
